# GIT Documentation.docx -- apply the two wording tweaks from the commit.
#
# 1) "By following these steps, you can set up and run each " becomes
#    "By following these steps, we should be able to set up and run each "
#    -- typed as three separate runs, with Word's auto-managed "_GoBack"
#    bookmark (last edit position) landing between "...able to " and
#    "set up...".
# 2) The old "_GoBack" location (inside "...indicate w|hether...") loses
#    the bookmark and the run split collapses back into a single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "By following these steps, you can set up and run each "
# ---------------------------------------------------------------------
$search = $d.Content
$found = $search.Find.Execute(
    "By following these steps, you can set up and run each ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sentStart = $search.Start

    # "you can " -> "we should be able to "
    $youCan = $d.Range($sentStart + 26, $sentStart + 34)
    $youCan.Text = "we should be able to "

    # Force a hard run boundary between "By following these steps, " and
    # "we should be able to " (matches the two separate <w:r> runs in the
    # target) by dropping a throwaway bookmark there and removing it again.
    $splitPoint = $d.Range($sentStart + 26, $sentStart + 26)
    $d.Bookmarks.Add("ZZZTempSplit", $splitPoint)
    $d.Bookmarks("ZZZTempSplit").Delete()

    # Word re-seats its hidden "_GoBack" bookmark at the last edit point --
    # that is the boundary right after "we should be able to ".
    $goBackAt = $sentStart + 26 + 21
    $goBackRange = $d.Range($goBackAt, $goBackAt)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}

# ---------------------------------------------------------------------
# Edit 2: merge "... indicate w" / "hether ..." back into one run and
# drop the bookmark that used to live on that boundary (handled above by
# re-adding "_GoBack" elsewhere, which moves it away from here).
# ---------------------------------------------------------------------
$search2 = $d.Content
$found2 = $search2.Find.Execute(
    "indicate whether they want to see their child's emotions",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $whole = $d.Range($search2.Start, $search2.End)
    # Round-trip through a placeholder so the engine re-segments the text
    # as a single run instead of leaving the old two-run split in place.
    $whole.Text = "ZZZPlaceholderForMerge"

    $search3 = $d.Content
    $found3 = $search3.Find.Execute(
        "ZZZPlaceholderForMerge",
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3) {
        $whole2 = $d.Range($search3.Start, $search3.End)
        $whole2.Text = "indicate whether they want to see their child's emotions"
    }
}
